# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.899.00"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "1.637.21"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'214.96"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "'0.0638"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.864.08"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.25"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.596.74"
$ws.Range("E14").Value = "  -3.82%  "

$ws.Range("D15").Value = "'0.544"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Value = "'62.88"
$ws.Range("E17").Value = "  -1.23%  "

$ws.Range("D18").Value = "25.914.27"
$ws.Range("E18").Value = "  -1.26%  "

$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").Value = "'193.05"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").Value = "'6.28"
$ws.Range("E23").Value = "  -1.10%  "

$ws.Range("D24").Value = "'143.91"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").Value = "'1.77"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("E27").Value = "  +1.62%  "

$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("D29").Value = "'15.51"
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").Value = "'3.31"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  -4.14%  "

$ws.Range("E35").Value = "  +0.96%  "

$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("D37").Value = "1.137.33"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("D40").Value = "'0.0158"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "'5.48"
$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("D43").Value = "'99.40"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").Value = "'0.800"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "1.773.72"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").Value = "0.0₆0115"
$ws.Range("E46").Value = "  +2.65%  "

$ws.Range("D47").Value = "'56.62"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("D48").Value = "'0.0533"
$ws.Range("E48").Value = "  +2.80%  "

$ws.Range("D49").Value = "'1.47"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = "  -0.57%  "

$ws.Range("D51").Value = "'0.414"
$ws.Range("E51").Value = "  -0.81%  "
